$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9598774909973145
$ws.Range("B1").Value = 2.076278448104858
$ws.Range("C1").Value = 5.085182189941406
$ws.Range("D1").Value = 1.898971557617188
$ws.Range("E1").Value = 1.324110746383667
